$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Notes" entries for the bumper / camera pins in column G
# (column G header already reads "Notes" at G1). Write in this order so
# new shared-string entries are interned left bumper, right bumper,
# right cam, back cam, left cam.
$ws.Range("G23").Value = "left bumper"
$ws.Range("G22").Value = "right bumper"
$ws.Range("G26").Value = "right cam"
$ws.Range("G27").Value = "back cam"
$ws.Range("G25").Value = "left cam"

# Match the active selection recorded in the saved workbook
$ws.Range("G26").Select()
